# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (copied from "2021-Q4" so it inherits the
#    same column layout/styles) positioned right before "总计".
# 2. Fill the new sheet with the 2022-Q1 fund-holding data.
# 3. Insert a new top data row into "总计" for the 2022-Q1 summary figures,
#    shifting the existing quarters down, and renumber the index column.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- 1) create the new sheet right before 总计 -----------------------------
$q4.Copy($total)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q1"

# the copy shuffled sheet positions, so re-resolve "总计" by name rather than
# trusting the (now stale) $total reference captured before the copy
$total = $wb.Worksheets.Item("总计")

# --- 2) populate 2022-Q1 ----------------------------------------------------
# fund code column keeps its leading zero -> force text before writing
$codeRng = $newSheet.Range("B2:B3")
$codeRng.NumberFormat = "@"
$newSheet.Range("B2").Value = "011685"
$newSheet.Range("B3").Value = "011686"
$codeRng.Style = "Normal"

$newSheet.Range("C2").Value = "创金合信先进装备股票A"
$newSheet.Range("C3").Value = "创金合信先进装备股票C"

# D:G are numeric-looking strings stored as text, same as the other quarters
$numStrRng = $newSheet.Range("D2:G3")
$numStrRng.NumberFormat = "@"
$newSheet.Range("D2").Value = "0.73"
$newSheet.Range("E2").Value = "92.01"
$newSheet.Range("F2").Value = "8.89"
$newSheet.Range("G2").Value = "0.0649"
$newSheet.Range("D3").Value = "0.17"
$newSheet.Range("E3").Value = "92.01"
$newSheet.Range("F3").Value = "8.89"
$newSheet.Range("G3").Value = "0.0151"
$numStrRng.Style = "Normal"

$newSheet.Range("H2").Value = 5
$newSheet.Range("H3").Value = 5

# --- 3) insert the new summary row into 总计 --------------------------------
$total.Rows(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.08
$total.Range("B2:D2").Style = "Normal"

# carry the header-row-style formatting onto the new index cell
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# renumber the index column for the rows pushed down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

Write-Output "2022-Q1 sheet added and 总计 updated"
